$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row: ID -> KundenNummer, Name -> Vorname, Stadt -> Ort (Alter stays)
$ws.Range("A1").Value = "KundenNummer"
$ws.Range("B1").Value = "Vorname"
$ws.Range("D1").Value = "Ort"

# Update Clara's age from 28 to 29
$ws.Range("C2").Value = 29

# Update selection to active cell E10
$ws.Range("E10").Select()
